# Auto-generated Excel COM-interop script
# Applies the data refresh described by the commit diff:
#  - numeric "want-to-go" counts (column F) bumped on many rows across all 4 sheets
#  - two 520-themed shows in the "Yan Chu" sheet marked as not-for-sale
#  - "All types" sheet: old row 11 event removed, rows 12-19 shift up one,
#    and a brand-new event is appended as the new row 19
#
# Note: every text Value assignment is prefixed with a leading apostrophe
# (PowerShell string literal "'...") -- this is the same trick Excels own
# UI uses to force literal-text storage and prevents the host from silently
# reinterpreting date-shaped strings (e.g. "2024-05-20") as real dates.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1538
$ws.Range("F5").Value = 8090
$ws.Range("F6").Value = 4937
$ws.Range("F7").Value = 7256
$ws.Range("F8").Value = 888
$ws.Range("F9").Value = 217
$ws.Range("F11").Value = 1206
$ws.Range("F13").Value = 216
$ws.Range("F14").Value = 587
$ws.Range("F16").Value = 57
$ws.Range("F17").Value = 258
$ws.Range("F20").Value = 1323
$ws.Range("F21").Value = 1285
$ws.Range("F24").Value = 1302
$ws.Range("F26").Value = 175
$ws.Range("F28").Value = 33
$ws.Range("F30").Value = 238
$ws.Range("F31").Value = 1041
$ws.Range("F33").Value = 21
$ws.Range("F34").Value = 171
$ws.Range("F36").Value = 48
$ws.Range("F38").Value = 587
$ws.Range("F40").Value = 104
$ws.Range("F42").Value = 129
$ws.Range("F44").Value = 1211
$ws.Range("F45").Value = 640
$ws.Range("F46").Value = 179

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 153
$ws.Range("F23").Value = 648
$ws.Range("F30").Value = 908
$ws.Range("F33").Value = 627
$ws.Range("F36").Value = 134

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 717
$ws.Range("F7").Value = 233
$ws.Range("F8").Value = 115
$ws.Range("F9").Value = 1860
$ws.Range("F10").Value = 2762

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1538
$ws.Range("F6").Value = 717
$ws.Range("F7").Value = 8090
$ws.Range("F8").Value = 233
$ws.Range("F9").Value = 4937
$ws.Range("F10").Value = 7256
$ws.Range("F20").Value = 216
$ws.Range("F21").Value = 587
$ws.Range("F22").Value = 258
$ws.Range("F23").Value = 1323
$ws.Range("F24").Value = 1285
$ws.Range("F26").Value = 1302
$ws.Range("F27").Value = 175
$ws.Range("F28").Value = 33
$ws.Range("F33").Value = 171
$ws.Range("F38").Value = 587
$ws.Range("F39").Value = 627
$ws.Range("F40").Value = 129
$ws.Range("F41").Value = 134
$ws.Range("F42").Value = 451
$ws.Range("F43").Value = 640
$ws.Range("F45").Value = 179

# ---- Sheet: 演出 : mark 520-themed shows as not-for-sale ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G10").Value = "'不可售"
$ws.Range("G11").Value = "'不可售"

# ---- Sheet: 全部类型 : rows 11-19 content refresh ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("B11").Value = "'2024-05-19"
$ws.Range("C11").Value = "'上海·恋与深空×恋与制作人only"
$ws.Range("D11").Value = "'顾村镇蕰川路6号 智慧湾科创园"
$ws.Range("E11").Value = "'2024.05.19 10:00-05.19 17:00"
$ws.Range("F11").Value = 888
$ws.Range("G11").Value = 70
$ws.Range("H11").Value = "'https://show.bilibili.com/platform/detail.html?id=83598"
$ws.Range("I11").Value = "'//i2.hdslb.com/bfs/openplatform/202403/jskKqUvJ1711165688442.jpeg"

$ws.Range("B12").Value = "'2024-05-20"
$ws.Range("C12").Value = "'上海·战双帕弥什 x HAPPY ZOO主题Cafe"
$ws.Range("D12").Value = "'广中路街道花园路128号 虹口德必运动LOFT"
$ws.Range("E12").Value = "'2024.05.20 00:00-06.03 23:59"
$ws.Range("F12").Value = 115
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = "'https://show.bilibili.com/platform/detail.html?id=85642"
$ws.Range("I12").Value = "'//i2.hdslb.com/bfs/openplatform/202405/ehdpyJQp1715739281505.png"

$ws.Range("B13").Value = "'2024-05-25"
$ws.Range("C13").Value = "'上海·Anime Market 同人展"
$ws.Range("D13").Value = "'漕宝路1688号 诺宝中心酒店"
$ws.Range("E13").Value = "'2024.05.25 08:00-05.26 20:00"
$ws.Range("F13").Value = 217
$ws.Range("G13").Value = 78
$ws.Range("H13").Value = "'https://show.bilibili.com/platform/detail.html?id=84737"
$ws.Range("I13").Value = "'//i1.hdslb.com/bfs/openplatform/202404/DyXYekek1713284815372.png"

$ws.Range("C14").Value = "'上海·EVANGELION× PrismLand · 新世纪福音战士官方授权主题店"
$ws.Range("D14").Value = "'南京东路830号第一百货商业中心B馆5楼(海底捞旁边) 第一百货商业中心"
$ws.Range("E14").Value = "'2024.05.25 00:00-07.22 23:59"
$ws.Range("F14").Value = 1860
$ws.Range("G14").Value = 20
$ws.Range("H14").Value = "'https://show.bilibili.com/platform/detail.html?id=85030"
$ws.Range("I14").Value = "'//i0.hdslb.com/bfs/openplatform/202404/K3kIpfaB1714445776157.jpeg"

$ws.Range("C15").Value = "'上海·「排球少年!!垃圾场决战 × animate cafe」"
$ws.Range("D15").Value = "'西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$ws.Range("E15").Value = "'2024.05.25 00:00-07.02 23:59"
$ws.Range("F15").Value = 2762
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = "'https://show.bilibili.com/platform/detail.html?id=85283"
$ws.Range("I15").Value = "'//i0.hdslb.com/bfs/openplatform/202405/vy2vecK11715162037223.jpeg"

$ws.Range("C16").Value = "'上海·你喜欢什么颜色呢？~黑泽朋世上海演唱会"
$ws.Range("D16").Value = "'中兴路1599号金融街融泰中心 蜚声上海PHASE LIVE HOUSE"
$ws.Range("E16").Value = "'2024.05.25 19:30-05.25 21:00"
$ws.Range("F16").Value = 214
$ws.Range("G16").Value = 580
$ws.Range("H16").Value = "'https://show.bilibili.com/platform/detail.html?id=83997"
$ws.Range("I16").Value = "'//i2.hdslb.com/bfs/openplatform/202404/0ias3gVf1714027450047.jpeg"

$ws.Range("C17").Value = "'上海·拉帮结派ONLY"
$ws.Range("D17").Value = "'海潮路133号B1 JUMP工坊"
$ws.Range("E17").Value = "'2024.05.25 14:00-05.25 19:00"
$ws.Range("F17").Value = 84
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = "'https://show.bilibili.com/platform/detail.html?id=85091"
$ws.Range("I17").Value = "'//i0.hdslb.com/bfs/openplatform/202405/Gqv3tfiB1714795562310.jpeg"

$ws.Range("C18").Value = "'上海·第六届Redamancy动漫游戏嘉年华"
$ws.Range("D18").Value = "'中山北路3300号4楼 上海环球港"
$ws.Range("E18").Value = "'2024.05.25 10:00-05.26 17:00"
$ws.Range("F18").Value = 1206
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = "'https://show.bilibili.com/platform/detail.html?id=84632"
$ws.Range("I18").Value = "'//i1.hdslb.com/bfs/openplatform/202404/im8C39eo1713190504331.png"

$ws.Range("B19").Value = "'2024-05-26"
$ws.Range("C19").Value = "'上海·五十岚隼士&小池亮介·2024见面会"
$ws.Range("D19").Value = "'长寿路街道万航渡后路19号 上海瓦肆文化传播有限公司"
$ws.Range("E19").Value = "'2024.05.26 11:00-05.26 15:30"
$ws.Range("F19").Value = 153
$ws.Range("G19").Value = 480
$ws.Range("H19").Value = "'https://show.bilibili.com/platform/detail.html?id=84615"
$ws.Range("I19").Value = "'//i0.hdslb.com/bfs/openplatform/202404/LwpFUbIJ1713857706981.jpeg"
